$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reformat the three summary rows whose status text changes from "Pending" to
#     "Not Recorded" so they use the same style as the other "Not Recorded" rows ---
# (style index 9 - fontId=2/fillId=5 - is retired in favor of style index 4 - fontId=2/fillId=4)
$styleSource = $ws.Range("A3:I3")
$styleSource.Copy()
$ws.Range("A89:I89").PasteSpecial(-4122)
$ws.Range("A127:I127").PasteSpecial(-4122)
$ws.Range("A167:I167").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Cell value / text updates ---
$ws.Range('G2').Value = 'servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range('G4').Value = 'hend_mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range('G5').Value = 'AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range('G6').Value = 'NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range("L7").Value = 45
$ws.Range("L8").Value = 0
$ws.Range('G9').Value = 'norhan.mohamed@med.asu.edu.eg, Ayaessam@med.asu.edu.eg'
$ws.Range('G10').Value = 'Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range('L10').Value = "'38.4%"
$ws.Range('G11').Value = 'Sarah.Mahdy@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Range('S16').Value = "'40.2%"
$ws.Range('G18').Value = 'aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("P18").Value = 5
$ws.Range("Q18").Value = 0
$ws.Range('G19').Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('S19').Value = "'48.3%"
$ws.Range('G20').Value = 'nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range("P20").Value = 7
$ws.Range("Q20").Value = 0
$ws.Range('S20').Value = "'44.5%"
$ws.Range('G21').Value = 'yasmin.m.senosy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range('G22').Value = 'naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("P22").Value = 6
$ws.Range("Q22").Value = 0
$ws.Range('G23').Value = 'yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Range('G24').Value = 'servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range('G25').Value = 'gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, backup@backdoor.com, manar.montaser@med.asu.edu.eg'
$ws.Range('G26').Value = 'hend_mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range('G27').Value = 'AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range('H27').Value = '142/217'
$ws.Range('G28').Value = 'NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Range('G32').Value = 'Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range('G33').Value = 'Sarah.Mahdy@med.asu.edu.eg, Alaa.A.Mostafa@med.asu.edu.eg'
$ws.Range('G40').Value = 'aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range('G41').Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G42').Value = 'nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range('G43').Value = 'yasmin.m.senosy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range('G44').Value = 'naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G45').Value = 'yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Range('G46').Value = 'hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range('G48').Value = 'heba@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G52').Value = 'mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
$ws.Range('G54').Value = 'yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range('G60').Value = 'nancy.abdelshafy@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg'
$ws.Range('G62').Value = 'wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range('G63').Value = 'ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg'
$ws.Range('G64').Value = 'wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range('G65').Value = 'nardine.alfonse@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G66').Value = 'marinasorial@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G68').Value = 'hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range('G70').Value = 'heba@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G74').Value = 'mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
$ws.Range('G75').Value = 'norhan.mohamed@med.asu.edu.eg, Ayaessam@med.asu.edu.eg'
$ws.Range('G76').Value = 'yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range('G81').Value = 'enas.omran@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range('G82').Value = 'nancy.abdelshafy@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg'
$ws.Range('G84').Value = 'wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range('G85').Value = 'ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg'
$ws.Range('G86').Value = 'wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range('G87').Value = 'nardine.alfonse@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G88').Value = 'marinasorial@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('I89').Value = 'Not Recorded'
$ws.Range('G92').Value = 'heba@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G93').Value = 'AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range('G96').Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range('G97').Value = 'norhan.mohamed@med.asu.edu.eg, Ayaessam@med.asu.edu.eg'
$ws.Range('G98').Value = 'yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range('G104').Value = 'nancy.abdelshafy@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg'
$ws.Range('G106').Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G107').Value = 'wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G108').Value = 'nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range('G110').Value = 'yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G111').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('H111').Value = '24/154'
$ws.Range('G114').Value = 'heba@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G115').Value = 'AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range('G118').Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range('G119').Value = 'norhan.mohamed@med.asu.edu.eg, Ayaessam@med.asu.edu.eg'
$ws.Range('G120').Value = 'yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range('G126').Value = 'nancy.abdelshafy@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg'
$ws.Range('H126').Value = '46/224'
$ws.Range('I127').Value = 'Not Recorded'
$ws.Range('G128').Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G129').Value = 'wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G130').Value = 'nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range('G131').Value = 'nardine.alfonse@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G132').Value = 'yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G133').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G134').Value = 'hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range('G137').Value = 'AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range('G141').Value = 'norhan.mohamed@med.asu.edu.eg, Ayaessam@med.asu.edu.eg'
$ws.Range('G142').Value = 'yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg'
$ws.Range('G150').Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G151').Value = 'wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G152').Value = 'Wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G153').Value = 'marinasorial@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G154').Value = 'naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G155').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G156').Value = 'alshimaa.atef@med.asu.edu.egm, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range('G159').Value = 'AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range('G163').Value = 'norhan.mohamed@med.asu.edu.eg, Ayaessam@med.asu.edu.eg'
$ws.Range('G164').Value = 'yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg'
$ws.Range('I167').Value = 'Not Recorded'
$ws.Range('G172').Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G173').Value = 'wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G174').Value = 'Wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G175').Value = 'marinasorial@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G176').Value = 'naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G177').Value = 'eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'

# --- Restore original (non quote-prefixed) formatting on the text-percentage cells ---
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("R16").Copy()
$ws.Range("S16").PasteSpecial(-4122)
$ws.Range("R19").Copy()
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("R20").Copy()
$ws.Range("S20").PasteSpecial(-4122)
$excel.CutCopyMode = 0
